$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 101
$ws.Range("B9").Value = "No more problems"
$ws.Range("C9").Value = "Database X"
$ws.Range("D9").Value = "Central Information Team"
$ws.Range("F9").Value = "Integer"
$ws.Range("G9").Value = "Decrease"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = "Division"
$ws.Range("J9").Value = "GGGG"
$ws.Range("K9").Value = 24

$ws.Range("A9").Select() | Out-Null
